$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "-2"
$ws.Range("B1").Value = "-1"
$ws.Range("A2").Value = "0.5"
$ws.Range("B2").Value = "1.5"
$ws.Range("A3").Value = "-1.5"
$ws.Range("B3").Value = "3.5"

$ws.Range("C1:D4").Clear()
$ws.Range("A4:B4").Clear()
